$wb = $excel.ActiveWorkbook

# Update the StaffMember sheet: replace staff member name and title in row 2
$staffSheet = $wb.Worksheets.Item("StaffMember")
$staffSheet.Range("B2").Value = "Executive Administrator"
$staffSheet.Range("A2").Value = "Nicole Bicho"

# Make StaffMember the active sheet/tab (it becomes tabSelected="1")
$staffSheet.Activate()
$staffSheet.Range("B13").Select()
